$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F19").Value = "18_hazards_to_humans_and_domestic_animals"
$ws.Range("F45").Value = "application instructions || env warning - species || env warning - water"
$ws.Range("F46").Value = "pollinator"
$ws.Range("F47").Value = "pollinator"
$ws.Range("F57").Value = "mixing"
$ws.Range("F68").Value = "application instructions"
$ws.Range("F69").Value = "application instructions"
$ws.Range("F71").Value = "use restrictions"
$ws.Range("F72").Value = "application instructions"
$ws.Range("F73").Value = "application instructions"
$ws.Range("F76").Value = "application instructions"
$ws.Range("F77").Value = "application instructions"
$ws.Range("F78").Value = "application instructions"
$ws.Range("F79").Value = "application instructions"
$ws.Range("F80").Value = "application instructions"
$ws.Range("F82").Value = "application instructions"
$ws.Range("F85").Value = "application instructions"
$ws.Range("F97").Value = "application instructions"
$ws.Range("F99").Value = "application instructions"
$ws.Range("F100").Value = "application instructions"
$ws.Range("F101").Value = "application instructions"
$ws.Range("F102").Value = "application instructions"
$ws.Range("F103").Value = "application instructions"
$ws.Range("F105").Value = "application instructions"
$ws.Range("F107").Value = "use restrictions"
$ws.Range("F108").Value = "use restrictions"
$ws.Range("F109").Value = "application instructions"
$ws.Range("F110").Value = "use restrictions"
$ws.Range("F116").Value = "application instructions"
$ws.Range("F117").Value = "use restrictions"
$ws.Range("F119").Value = "mixing"
$ws.Range("F155").Value = "154_pesticide_storage"
